# New weekly price record for "Jengibre" (Mercado Mayorista Lo Valledor de
# Santiago) is inserted as a new data row right after the existing row for
# 2021-05-01 (row 37, date 44277) and before the former row 38 (date 44333).
# Inserting a whole row there shifts every following record down by one,
# which reproduces the rest of the diff (each subsequent row now holds the
# data that used to belong to the row above it), and grows the used range
# from A1:R53 to A1:R54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 38..53 down to 39..54, leaving a fresh blank row 38.
$ws.Rows(38).Insert()

# Populate the newly inserted row 38 with the new weekly record.
$ws.Cells.Item(38, 1).Value = 6
$ws.Cells.Item(38, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(38, 3).Value = "Metropolitana"
$ws.Cells.Item(38, 4).Value = 44508
$ws.Cells.Item(38, 5).Value = 13
$ws.Cells.Item(38, 6).Value = 100114007
$ws.Cells.Item(38, 7).Value = "Jengibre"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 400
$ws.Cells.Item(38, 11).Value = 13000
$ws.Cells.Item(38, 12).Value = 15000
$ws.Cells.Item(38, 13).Value = 13850
$ws.Cells.Item(38, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(38, 15).Value = "Perú"
$ws.Cells.Item(38, 16).Value = 1065
$ws.Cells.Item(38, 17).Value = 13
$ws.Cells.Item(38, 18).Value = "Hortaliza"
